$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the custom trait formula labels to match normal trait formatting
$ws.Range("B1").Value = "first_trait = H3N4 + H3N4F1 + H4N4 + H4N4F1"
$ws.Range("B2").Value = "second_trait = (0.5 * H3N4 + H4N4) / (H3N4F1 + H4N4F1)"

# Update selection to mirror the author's final cursor position
$ws.Range("D9").Select()
